# Generate Report for Handoff
#
# Refresh the localization-status report: the handoff batch that was
# previously pending now carries the "ht" (handoff type) priority, and
# the handoff timestamps recorded for that batch move forward to the
# moment the new xliff files were actually generated.
#
# Affected rows (Source File Name): 08cf076b-*, 1a57c63d-*, 6b6a5385-*,
# c00e38e5-*, ea649a8c-*, f9059e69-* -> table rows 7, 8, 11, 12, 13, 14
# on both the "zh-cn" and "de-de" worksheets, plus the corresponding
# "Latest HO Xliff Generate Date" column on the "Overview" worksheet.

$wb = $excel.ActiveWorkbook

$rows = @("7", "8", "11", "12", "13", "14")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
# mirrors the de-de handoff timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in $rows) {
    $wsOverview.Range("G$row").Value = "2016-08-16 04:22:53"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($row in $rows) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-08-16 04:22:48"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($row in $rows) {
    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-08-16 04:22:53"
}
